$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 110, shifting rows 110:184 down to 111:185.
$ws.Rows("110:110").Insert()

# Populate the newly inserted row 110 with the new record's data.
$ws.Range("A110").Value = 5
$ws.Range("B110").Value = "Macroferia Regional de Talca"
$ws.Range("C110").Value = "Maule"
$ws.Range("D110").Value = 44596
$ws.Range("D110").NumberFormat = $ws.Range("D111").NumberFormat()
$ws.Range("E110").Value = 7
$ws.Range("F110").Value = 100112024
$ws.Range("G110").Value = "Choclo"
$ws.Range("H110").Value = "Choclero"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 40000
$ws.Range("K110").Value = 100
$ws.Range("L110").Value = 100
$ws.Range("M110").Value = 100
$ws.Range("N110").Value = "$/unidad"
$ws.Range("O110").Value = "Región del Maule"
$ws.Range("P110").Value = 100
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"
